# Week 15 simulations: a new QB (N.Mullens) was added to the Browns roster.
# On the "QB" sheet, insert a new row above B.Mayfield's row for N.Mullens
# with all his Week 15 stat columns at 0, then leave the view focused on
# that sheet (matching the post-edit state captured in the workbook).

$wb = $excel.ActiveWorkbook
$qb = $wb.Worksheets.Item("QB")

# Push the existing players (B.Mayfield, C.Keenum) down one row and make
# room for the new player's stat line right under the header row.
$qb.Rows.Item(2).Insert() | Out-Null

$qb.Range("A2").Value = "N.Mullens"
$qb.Range("B2:L2").Value = 0

# Reflect the updated selection/active-sheet state left behind by the edit.
$qb.Select() | Out-Null
$qb.Range("M5").Select() | Out-Null

Write-Output "Added N.Mullens to QB sheet (Week 15)"
